# Fruta / hortaliza, semanal
# Insert two new weekly rows of data just before the existing row 221
# (Femacal de La Calera - Frutilla), shifting the remaining historical
# rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above what is currently row 221. This pushes the
# old rows 221-307 down to 223-309 and keeps their formatting/values intact.
$ws.Rows("221:222").Insert()

# --- New row 221 ---
$ws.Cells.Item(221, 1).Value = 3
$ws.Cells.Item(221, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(221, 3).Value = "Coquimbo"
$ws.Cells.Item(221, 4).Value = 44704
$ws.Cells.Item(221, 5).Value = 5
$ws.Cells.Item(221, 6).Value = "Fruta"
$ws.Cells.Item(221, 7).Value = 100101
$ws.Cells.Item(221, 8).Value = "Berries"
$ws.Cells.Item(221, 9).Value = 100112025
$ws.Cells.Item(221, 10).Value = "Frutilla"
$ws.Cells.Item(221, 11).Value = "Sin especificar"
$ws.Cells.Item(221, 12).Value = "Especial"
$ws.Cells.Item(221, 13).Value = 45
$ws.Cells.Item(221, 14).Value = 10000
$ws.Cells.Item(221, 15).Value = 10000
$ws.Cells.Item(221, 16).Value = 10000
$ws.Cells.Item(221, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(221, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(221, 19).Value = 1429
$ws.Cells.Item(221, 20).Value = 7

# --- New row 222 ---
$ws.Cells.Item(222, 1).Value = 3
$ws.Cells.Item(222, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(222, 3).Value = "Coquimbo"
$ws.Cells.Item(222, 4).Value = 44704
$ws.Cells.Item(222, 5).Value = 5
$ws.Cells.Item(222, 6).Value = "Fruta"
$ws.Cells.Item(222, 7).Value = 100101
$ws.Cells.Item(222, 8).Value = "Berries"
$ws.Cells.Item(222, 9).Value = 100112025
$ws.Cells.Item(222, 10).Value = "Frutilla"
$ws.Cells.Item(222, 11).Value = "Sin especificar"
$ws.Cells.Item(222, 12).Value = "Primera"
$ws.Cells.Item(222, 13).Value = 40
$ws.Cells.Item(222, 14).Value = 7000
$ws.Cells.Item(222, 15).Value = 7000
$ws.Cells.Item(222, 16).Value = 7000
$ws.Cells.Item(222, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(222, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(222, 19).Value = 1000
$ws.Cells.Item(222, 20).Value = 7
